## backup before deadline cleaning
## Rebuilds row 1 (headers) and row 2 (values) of Sheet1 so that the
## sheet now reports the full q-quantile set (71 columns, A:BS) instead
## of the previous 45-column subset (A:AS).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading "'" forces the text/quote-prefix style (s="1") that the header
# row already uses, while the shared-string value stored is the plain
# "qN" label (the apostrophe is just an entry marker, not part of the text).
$headers = @("'q0","'q1","'q2","'q3","'q6","'q7","'q8","'q9","'q10","'q11","'q16","'q17","'q18","'q19","'q21","'q22","'q23","'q24","'q25","'q26","'q27","'q28","'q29","'q30","'q31","'q32","'q33","'q34","'q35","'q36","'q37","'q38","'q39","'q40","'q41","'q42","'q43","'q44","'q45","'q46","'q47","'q48","'q49","'q50","'q51","'q52","'q53","'q54","'q55","'q56","'q57","'q58","'q59","'q60","'q61","'q62","'q63","'q64","'q65","'q66","'q67","'q68","'q69","'q70","'q74","'q75","'q76","'q77","'q78","'q79","'q86")

$values = @(0.346,0.583,0.367,0.5720000000000001,0.206,0.099,0.332,0.501,0.022,0.587,0.126,0.032,0.206,0.399,0.335,0.324,0.399,0.8350000000000001,0.28,0.071,0.334,0.399,1.3840000000000001,0.888,0.251,0.032,0.542,1.524,1.2930000000000001,0.7250000000000001,0.5800000000000001,0.7050000000000001,1.3860000000000001,1.2640000000000002,0.803,0.44499999999999995,1.826,0.793,0.6010000000000001,0.547,1.1400000000000001,0.446,0.28300000000000003,0.4,0.399,0.5509999999999999,0.27,0.187,0.333,0.9410000000000001,0.8280000000000001,1.099,0.5980000000000001,0.734,1.205,0.538,0.703,0.16699999999999998,0.138,0.01,0.087,0.092,0.399,0.798,0.399,0.39,0.113,0.399,0.353,0.151,0.399)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
    $ws.Cells.Item(2, $col).Value = $values[$i]
}
